$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of mod-count data for 2025/12/12.
$ws.Range("A33").NumberFormat = "@"
$ws.Range("A33").Value = "2025/12/12"
$ws.Range("B33").Value = "逃离鸭科夫"
$ws.Range("C33").Value = 1356

# Match the formatting (centered alignment style) used by the other data rows.
$ws.Range("A32:C32").Copy()
$ws.Range("A33:C33").PasteSpecial(-4122)
